$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 255.25
$ws.Range("I4").Value = 215
$ws.Range("J4").Value = 376
$ws.Range("K4").Value = 215
$ws.Range("L4").Value = 376
$ws.Range("M4").Value = -101
$ws.Range("N4").Value = -604

$ws.Range("H33").Value = 35722300
$ws.Range("I33").Value = 71429976
$ws.Range("J33").Value = 14621.143
$ws.Range("K33").Value = 71429976
$ws.Range("L33").Value = 14621.143
$ws.Range("M33").Value = -71429747
$ws.Range("N33").Value = -15079.143

$ws.Range("H86").Value = 2913.4
$ws.Range("I86").Value = 3237.875
$ws.Range("J86").Value = 2542.5715
$ws.Range("K86").Value = 3237.875
$ws.Range("L86").Value = 2542.5715
$ws.Range("M86").Value = -2114.875
$ws.Range("N86").Value = -4788.5715

$ws.Range("H89").Value = 2913.4
$ws.Range("I89").Value = 3237.875
$ws.Range("J89").Value = 2542.5715
$ws.Range("K89").Value = 16189.375
$ws.Range("L89").Value = 12712.8575
$ws.Range("M89").Value = -10573.375
$ws.Range("N89").Value = -23944.8575

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 754.65625
$ws.Range("I74").Value = 628.8077
$ws.Range("J74").Value = 1300
$ws.Range("K74").Value = 628.8077
$ws.Range("L74").Value = 1300
$ws.Range("M74").Value = 245.1923
$ws.Range("N74").Value = -3048

$ws.Range("H77").Value = 754.65625
$ws.Range("I77").Value = 628.8077
$ws.Range("J77").Value = 1300
$ws.Range("K77").Value = 3144.0385
$ws.Range("L77").Value = 6500
$ws.Range("M77").Value = 1223.9615
$ws.Range("N77").Value = -15236

$ws.Range("H112").Value = 17591.334
$ws.Range("J112").Value = 17591.334
$ws.Range("L112").Value = 17591.334
$ws.Range("N112").Value = -20545.334

$ws.Range("H122").Value = 2082.125
$ws.Range("I122").Value = 1454.2609
$ws.Range("K122").Value = 4362.7827
$ws.Range("M122").Value = -1912.7827

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1255549.5
$ws.Range("I4").Value = 1255.7142
$ws.Range("J4").Value = 2231111.2
$ws.Range("K4").Value = 1255.7142
$ws.Range("L4").Value = 2231111.2
$ws.Range("M4").Value = -1143.7142
$ws.Range("N4").Value = -2231335.2

$ws.Range("H16").Value = 5094.5625
$ws.Range("I16").Value = 2866.6667
$ws.Range("J16").Value = 6431.3
$ws.Range("K16").Value = 2866.6667
$ws.Range("L16").Value = 6431.3
$ws.Range("M16").Value = -2579.6667
$ws.Range("N16").Value = -7005.3

$ws.Range("H31").Value = 1988.04
$ws.Range("I31").Value = 1499.75
$ws.Range("J31").Value = 5568.8335
$ws.Range("K31").Value = 1499.75
$ws.Range("L31").Value = 5568.8335
$ws.Range("M31").Value = -1204.75
$ws.Range("N31").Value = -6158.8335

$ws.Range("H34").Value = 1988.04
$ws.Range("I34").Value = 1499.75
$ws.Range("J34").Value = 5568.8335
$ws.Range("K34").Value = 1499.75
$ws.Range("L34").Value = 5568.8335
$ws.Range("M34").Value = -1297.75
$ws.Range("N34").Value = -5972.8335

$ws.Range("H113").Value = 5094.5625
$ws.Range("I113").Value = 2866.6667
$ws.Range("J113").Value = 6431.3
$ws.Range("K113").Value = 2866.6667
$ws.Range("L113").Value = 6431.3
$ws.Range("M113").Value = -696.6667000000002
$ws.Range("N113").Value = -10771.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 737.8570999999999
$ws.Range("I5").Value = 603.55554
$ws.Range("J5").Value = 979.6
$ws.Range("K5").Value = 1810.66662
$ws.Range("L5").Value = 2938.8
$ws.Range("M5").Value = -1698.66662
$ws.Range("N5").Value = -3162.8

$ws.Range("H12").Value = 38.058823
$ws.Range("J12").Value = 48.083332
$ws.Range("L12").Value = 144.249996
$ws.Range("N12").Value = -490.249996

$ws.Range("H49").Value = 3375
$ws.Range("J49").Value = 3400
$ws.Range("L49").Value = 10200
$ws.Range("N49").Value = -10512

$ws.Range("H92").Value = 1112.5
$ws.Range("I92").Value = 550
$ws.Range("J92").Value = 1675
$ws.Range("K92").Value = 1650
$ws.Range("L92").Value = 5025
$ws.Range("M92").Value = -402
$ws.Range("N92").Value = -7521

$ws.Range("H122").Value = 2273178.8
$ws.Range("J122").Value = 7143558
$ws.Range("L122").Value = 64292022
$ws.Range("N122").Value = -64296922

$ws.Range("H131").Value = 869.4625
$ws.Range("I131").Value = 292
$ws.Range("J131").Value = 951.95715
$ws.Range("K131").Value = 876
$ws.Range("L131").Value = 2855.87145
$ws.Range("M131").Value = 4164
$ws.Range("N131").Value = -12935.87145

$ws.Range("H135").Value = 737.8570999999999
$ws.Range("I135").Value = 603.55554
$ws.Range("J135").Value = 979.6
$ws.Range("K135").Value = 5431.99986
$ws.Range("L135").Value = 8816.4
$ws.Range("M135").Value = -2896.99986
$ws.Range("N135").Value = -13886.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2516.5
$ws.Range("J80").Value = 2833
$ws.Range("L80").Value = 2833
$ws.Range("N80").Value = -4829

$ws.Range("H83").Value = 2516.5
$ws.Range("J83").Value = 2833
$ws.Range("L83").Value = 14165
$ws.Range("N83").Value = -24149

$ws.Range("H113").Value = 2600.8462
$ws.Range("I113").Value = 3218.5
$ws.Range("J113").Value = 2071.4285
$ws.Range("K113").Value = 3218.5
$ws.Range("L113").Value = 2071.4285
$ws.Range("M113").Value = -1048.5
$ws.Range("N113").Value = -6411.4285

$ws.Range("H122").Value = 2484.5881
$ws.Range("I122").Value = 2341.6
$ws.Range("K122").Value = 7024.799999999999
$ws.Range("M122").Value = -4574.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2618.4546
$ws.Range("I40").Value = 2387.875
$ws.Range("J40").Value = 3233.3333
$ws.Range("K40").Value = 2387.875
$ws.Range("L40").Value = 3233.3333
$ws.Range("M40").Value = -2251.875
$ws.Range("N40").Value = -3505.3333

$ws.Range("H122").Value = 2937.111
$ws.Range("I122").Value = 2900.0908
$ws.Range("J122").Value = 3100
$ws.Range("K122").Value = 8700.2724
$ws.Range("L122").Value = 9300
$ws.Range("M122").Value = -6250.2724
$ws.Range("N122").Value = -14200

$ws.Range("H132").Value = 2169.7856
$ws.Range("I132").Value = 1149
$ws.Range("J132").Value = 3828.5625
$ws.Range("K132").Value = 3447
$ws.Range("L132").Value = 11485.6875
$ws.Range("M132").Value = -917
$ws.Range("N132").Value = -16545.6875

$ws.Range("H136").Value = 1958.5
$ws.Range("I136").Value = 1011.75
$ws.Range("J136").Value = 3473.3
$ws.Range("K136").Value = 3035.25
$ws.Range("L136").Value = 10419.9
$ws.Range("M136").Value = -485.25
$ws.Range("N136").Value = -15519.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 22000
$ws.Range("J34").Value = 22000
$ws.Range("L34").Value = 22000
$ws.Range("N34").Value = -22406

$ws.Range("H136").Value = 2910.78
$ws.Range("I136").Value = 750.1852
$ws.Range("J136").Value = 5447.1304
$ws.Range("K136").Value = 2250.5556
$ws.Range("L136").Value = 16341.3912
$ws.Range("M136").Value = 299.4443999999999
$ws.Range("N136").Value = -21441.3912
